$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the title heading ("Play Cafelito Slot for Free - Review 2021").
# ------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description*") {
        $metaRange = $d.Range($p.Range.Start, $p.Range.End)
        $metaRange.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Cafelito Slot for Free - Review 2021"
#    right before the final ("Create a feature image ...") paragraph.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

# Insert the new paragraph plus an empty placeholder paragraph (needed so
# the runs of the following paragraph are not merged into this insertion).
$xmlFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cafelito Slot for Free - Review 2021</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>'
$insertPoint.InsertXML($xmlFragment) | Out-Null

# Remove the placeholder empty paragraph that was introduced purely to keep
# the paragraph boundary from merging with the next paragraph's runs.
$placeholderPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$placeholderRange = $d.Range($placeholderPara.Range.Start, $placeholderPara.Range.End)
if (($placeholderRange.End - $placeholderRange.Start) -eq 1) {
    $placeholderRange.Delete()
}

# ------------------------------------------------------------------
# 3) Update the text of the final (italic) paragraph from the old
#    "Create a feature image ..." prompt text to the meta-description text.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Create a feature image for Cafelito featuring a happy Maya warrior with glasses in a cartoon style. The background should have a coffee shop theme with images of coffee beans, cups, and machines. The Maya warrior should be holding a cup of coffee with a smile on their face. They should be wearing a colorful outfit with traditional Maya patterns, and their hair should be decorated with coffee beans and flowers. The image should be bright and vibrant, capturing the fun and excitement of playing the Cafelito slot game.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Read our review of Cafelito slot game and play for free. Discover the graphics, gameplay, bonuses, RTP value, and jackpots. Start playing now.", 2) | Out-Null
